$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "total"
$ws.Range("A5").Value = "urbana"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.2507872738196449
$ws.Range("D5").Value = 0.2484603367430894
$ws.Range("E5").Value = 0.7515164003538076
$ws.Range("F5").Value = 0.8169748258204544
$ws.Range("G5").Value = 1.940780142856299
$ws.Range("H5").Value = 0.9247170617093116
$ws.Range("A6").Value = "rural"
$ws.Range("C6").Value = 0.4173161102117761
$ws.Range("D6").Value = 0.4717473170636254
$ws.Range("E6").Value = 2.676333028325332
$ws.Range("F6").Value = 1.746744600816121
$ws.Range("G6").Value = 5.93810535821779
$ws.Range("H6").Value = 1.839644114289289
$ws.Range("A7").Value = "norte"
$ws.Range("C7").Value = 0.5735831757557439
$ws.Range("D7").Value = 0.6201527215424948
$ws.Range("E7").Value = 2.569956751788632
$ws.Range("F7").Value = 0.8076751940046523
$ws.Range("G7").Value = 5.400183655818156
$ws.Range("H7").Value = 0.8171824200618959
$ws.Range("A8").Value = "rondônia"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 1.787574207147729
$ws.Range("D8").Value = 1.845383093636107
$ws.Range("E8").Value = 3.280111635971241
$ws.Range("F8").Value = 1.917465776386419
$ws.Range("G8").Value = 9.152921625369048
$ws.Range("H8").Value = 2.166024832035191
$ws.Range("A9").Value = "acre"
$ws.Range("C9").Value = 1.499729312334564
$ws.Range("D9").Value = 1.731236975688517
$ws.Range("E9").Value = 8.360121698142493
$ws.Range("F9").Value = 2.980773848900957
$ws.Range("G9").Value = 21.17158299087262
$ws.Range("H9").Value = 2.299304691205551
$ws.Range("A10").Value = "amazonas"
$ws.Range("C10").Value = 1.172518521870318
$ws.Range("D10").Value = 1.194483463631898
$ws.Range("E10").Value = 6.730878941768616
$ws.Range("F10").Value = 1.612754390746685
$ws.Range("G10").Value = 24.16468413765944
$ws.Range("H10").Value = 1.490412236449429
$ws.Range("A11").Value = "roraima"
$ws.Range("C11").Value = 2.190731659551413
$ws.Range("D11").Value = 2.286381017219059
$ws.Range("E11").Value = 7.297101759456418
$ws.Range("F11").Value = 2.699326120370863
$ws.Range("G11").Value = 17.49024265404395
$ws.Range("H11").Value = 2.949951524630713
$ws.Range("A12").Value = "pará"
$ws.Range("C12").Value = 1.100868528515508
$ws.Range("D12").Value = 1.229038321734609
$ws.Range("E12").Value = 4.181318688618839
$ws.Range("F12").Value = 1.184341224243167
$ws.Range("G12").Value = 8.698134797846837
$ws.Range("H12").Value = 1.228725908267395
$ws.Range("A13").Value = "amapá"
$ws.Range("C13").Value = 2.296499906869757
$ws.Range("D13").Value = 2.436315927189144
$ws.Range("E13").Value = 11.08168071728537
$ws.Range("F13").Value = 4.364888746986128
$ws.Range("G13").Value = 27.02937796853413
$ws.Range("H13").Value = 3.838384746143414
$ws.Range("A14").Value = "tocantins"
$ws.Range("C14").Value = 0.8747345659301526
$ws.Range("D14").Value = 1.011497318577143
$ws.Range("E14").Value = 6.99081863238905
$ws.Range("F14").Value = 2.35371806581774
$ws.Range("G14").Value = 13.45456944742002
$ws.Range("H14").Value = 2.523408652592504
$ws.Range("A15").Value = "nordeste"
$ws.Range("C15").Value = 0.345036926675498
$ws.Range("D15").Value = 0.3531658569098057
$ws.Range("E15").Value = 2.389546029110432
$ws.Range("F15").Value = 0.9285988828769437
$ws.Range("G15").Value = 4.327064832285274
$ws.Range("H15").Value = 1.022481618902623
$ws.Range("A16").Value = "maranhão"
$ws.Range("C16").Value = 1.123903642695248
$ws.Range("D16").Value = 1.139879529585017
$ws.Range("E16").Value = 11.45538443957088
$ws.Range("F16").Value = 3.823721885785956
$ws.Range("G16").Value = 20.78756985864832
$ws.Range("H16").Value = 4.773019298692292
$ws.Range("A17").Value = "piauí"
$ws.Range("C17").Value = 1.330320961449737
$ws.Range("D17").Value = 1.328669414632295
$ws.Range("E17").Value = 9.665932839855074
$ws.Range("F17").Value = 3.063843742694349
$ws.Range("G17").Value = 29.32862853401107
$ws.Range("H17").Value = 2.611761322175687
$ws.Range("A18").Value = "ceará"
$ws.Range("C18").Value = 0.9262219973629479
$ws.Range("D18").Value = 0.9664486809695677
$ws.Range("E18").Value = 4.023148016102164
$ws.Range("F18").Value = 1.583714020253246
$ws.Range("G18").Value = 12.16154725360941
$ws.Range("H18").Value = 1.78882064903994
$ws.Range("A19").Value = "rio grande do norte"
$ws.Range("C19").Value = 0.9709016469290312
$ws.Range("D19").Value = 1.045334713374307
$ws.Range("E19").Value = 3.331784778822068
$ws.Range("F19").Value = 1.712963058890707
$ws.Range("G19").Value = 14.11365520358395
$ws.Range("H19").Value = 1.814214274943321
$ws.Range("A20").Value = "paraíba"
$ws.Range("C20").Value = 1.213185675589167
$ws.Range("D20").Value = 1.222953676197773
$ws.Range("E20").Value = 5.886463061944519
$ws.Range("F20").Value = 2.904989570309735
$ws.Range("G20").Value = 11.72960131123549
$ws.Range("H20").Value = 2.548242870724642
$ws.Range("A21").Value = "pernambuco"
$ws.Range("C21").Value = 0.8889787308692259
$ws.Range("D21").Value = 0.8669790398943953
$ws.Range("E21").Value = 3.71438129064304
$ws.Range("F21").Value = 1.861732782959334
$ws.Range("G21").Value = 9.363769405098003
$ws.Range("H21").Value = 2.064901127434388
$ws.Range("A22").Value = "alagoas"
$ws.Range("C22").Value = 1.067905635753283
$ws.Range("D22").Value = 1.052688204229557
$ws.Range("E22").Value = 14.6777990374004
$ws.Range("F22").Value = 4.847704027268084
$ws.Range("G22").Value = 23.78119527638065
$ws.Range("H22").Value = 5.165723610661686
$ws.Range("A23").Value = "sergipe"
$ws.Range("C23").Value = 1.842239894430963
$ws.Range("D23").Value = 1.843704314060877
$ws.Range("E23").Value = 6.69942499725019
$ws.Range("F23").Value = 2.654282299821404
$ws.Range("G23").Value = 18.31989754188432
$ws.Range("H23").Value = 2.323168578523941
$ws.Range("A24").Value = "bahia"
$ws.Range("C24").Value = 0.6728705371677834
$ws.Range("D24").Value = 0.7121572906296588
$ws.Range("E24").Value = 5.781730474822889
$ws.Range("F24").Value = 1.906466360518141
$ws.Range("G24").Value = 4.962785242991031
$ws.Range("H24").Value = 2.139183204710935
$ws.Range("A25").Value = "sudeste"
$ws.Range("C25").Value = 0.3785872225912779
$ws.Range("D25").Value = 0.3797503718171671
$ws.Range("E25").Value = 0.9817592655071447
$ws.Range("F25").Value = 1.310917780733238
$ws.Range("G25").Value = 2.835189599922381
$ws.Range("H25").Value = 1.52910793501635
$ws.Range("A26").Value = "minas gerais"
$ws.Range("C26").Value = 0.6609235877442553
$ws.Range("D26").Value = 0.6683793116574943
$ws.Range("E26").Value = 2.098519024071414
$ws.Range("F26").Value = 1.715542775836872
$ws.Range("G26").Value = 4.987158154972104
$ws.Range("H26").Value = 2.061489444748408
$ws.Range("A27").Value = "espírito santo"
$ws.Range("C27").Value = 1.805987988009141
$ws.Range("D27").Value = 1.890556741606611
$ws.Range("E27").Value = 6.500471262713051
$ws.Range("F27").Value = 5.094830251506438
$ws.Range("G27").Value = 13.73722224131833
$ws.Range("H27").Value = 4.936899525025382
$ws.Range("A28").Value = "rio de janeiro"
$ws.Range("C28").Value = 0.7743814032865057
$ws.Range("D28").Value = 0.7211946977886223
$ws.Range("E28").Value = 2.347915540304146
$ws.Range("F28").Value = 2.84781134607361
$ws.Range("G28").Value = 5.032334825523151
$ws.Range("H28").Value = 3.374686302400915
$ws.Range("A29").Value = "são paulo"
$ws.Range("C29").Value = 0.5733360619825494
$ws.Range("D29").Value = 0.5841402883519828
$ws.Range("E29").Value = 1.284134114785361
$ws.Range("F29").Value = 2.501275180161155
$ws.Range("G29").Value = 4.783683700293635
$ws.Range("H29").Value = 2.915105094593925
$ws.Range("A30").Value = "sul"
$ws.Range("C30").Value = 0.5093331640658189
$ws.Range("D30").Value = 0.5146519024045204
$ws.Range("E30").Value = 0.8558635976575422
$ws.Range("F30").Value = 3.147285091904333
$ws.Range("G30").Value = 5.283490065603138
$ws.Range("H30").Value = 3.464380743312173
$ws.Range("A31").Value = "paraná"
$ws.Range("C31").Value = 0.9920714210558192
$ws.Range("D31").Value = 1.002223587435989
$ws.Range("E31").Value = 1.777686182806633
$ws.Range("F31").Value = 4.100161680489536
$ws.Range("G31").Value = 9.148034368236482
$ws.Range("H31").Value = 4.42343458880582
$ws.Range("A32").Value = "santa catarina"
$ws.Range("C32").Value = 0.8308259620595784
$ws.Range("D32").Value = 0.8589335070446785
$ws.Range("E32").Value = 1.525734433561198
$ws.Range("F32").Value = 9.68726404712152
$ws.Range("G32").Value = 14.87338282674572
$ws.Range("H32").Value = 10.587113692343
$ws.Range("A33").Value = "rio grande do sul"
$ws.Range("C33").Value = 0.7854226609867928
$ws.Range("D33").Value = 0.7824129561702254
$ws.Range("E33").Value = 1.018353154528368
$ws.Range("F33").Value = 4.445270606584808
$ws.Range("G33").Value = 6.840934290718899
$ws.Range("H33").Value = 5.042361736176879
$ws.Range("A34").Value = "centro-oeste"
$ws.Range("C34").Value = 0.750520663841569
$ws.Range("D34").Value = 0.7804996092149774
$ws.Range("E34").Value = 2.368011155705495
$ws.Range("F34").Value = 1.722070590827341
$ws.Range("G34").Value = 4.760347776508738
$ws.Range("H34").Value = 1.822946403447939
$ws.Range("A35").Value = "mato grosso do sul"
$ws.Range("C35").Value = 1.487495043654953
$ws.Range("D35").Value = 1.563410657302831
$ws.Range("E35").Value = 3.798906638933747
$ws.Range("F35").Value = 3.827483236914496
$ws.Range("G35").Value = 11.94853646827792
$ws.Range("H35").Value = 4.038313037684235
$ws.Range("A36").Value = "mato grosso"
$ws.Range("C36").Value = 1.623718621213575
$ws.Range("D36").Value = 1.788795126644249
$ws.Range("E36").Value = 8.139903460204051
$ws.Range("F36").Value = 4.880184218354154
$ws.Range("G36").Value = 8.686973923858922
$ws.Range("H36").Value = 5.335986282966507
$ws.Range("A37").Value = "goiás"
$ws.Range("C37").Value = 1.265575303369433
$ws.Range("D37").Value = 1.312081979881054
$ws.Range("E37").Value = 3.347077237440443
$ws.Range("F37").Value = 2.22351756096297
$ws.Range("G37").Value = 7.967463943914436
$ws.Range("H37").Value = 2.362897401430239
$ws.Range("A38").Value = "distrito federal"
$ws.Range("C38").Value = 1.230162831797159
$ws.Range("D38").Value = 1.164996213384112
$ws.Range("E38").Value = 3.105697902105879
$ws.Range("F38").Value = 2.444684169609276
$ws.Range("G38").Value = 8.529075094521358
$ws.Range("H38").Value = 2.689108132926339

$ws.Range("A39:A40").EntireRow.Delete()
Write-Output "edit applied"